$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.283.70"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.872.78"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7108"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.61"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +1.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07730"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.05"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08458"
$ws.Range("E11").Value = "  +2.50%  "

$ws.Range("D12").Value = "1.871.42"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.201"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7110"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "29.288.76"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008311"
$ws.Range("E17").Value = "  +6.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.987"
$ws.Range("E18").Value = "  +2.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.39"
$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.122.87"
$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.809"
$ws.Range("E23").Value = "  -1.97%  "

$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1606"
$ws.Range("E25").Value = "  +2.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.14"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.014"
$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  +1.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.514"
$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.404"
$ws.Range("E30").Value = "  +1.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.320"
$ws.Range("E31").Value = "  +5.56%  "

$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05258"
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.919"
$ws.Range("E34").Value = "  +1.07%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7440"
$ws.Range("E36").Value = "  +2.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.712"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").Value = "1.167.59"
$ws.Range("E40").Value = "  +2.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.357"
$ws.Range("E41").Value = "  +4.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.85"
$ws.Range("E42").Value = "  +0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8848"
$ws.Range("E43").Value = "  -1.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.60"
$ws.Range("E44").Value = "  +4.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "2.018.52"
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.810"
$ws.Range("E47").Value = "  +2.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5198"
$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000121"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.370"
$ws.Range("E50").Value = "  +0.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4299"
$ws.Range("E51").Value = "  +1.41%  "

